$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.386665105819702
$ws.Range("B1").Value = 1.578549027442932
$ws.Range("C1").Value = 6.932479858398438
$ws.Range("D1").Value = 1.925748705863953
$ws.Range("E1").Value = 0.8726926445960999
